$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: AD1 (Wins), AE1 (Losses), AF1 (Ties)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold, centered, bordered style) from an
# existing header cell (AC1) onto the three new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill season-record columns for every data row (2-58) with the team's
# Wins / Losses / Ties totals.
$ws.Range("AD2:AD58").Value = 88
$ws.Range("AE2:AE58").Value = 74
$ws.Range("AF2:AF58").Value = 0
